$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full corrected dataset (fixes overwrite bug that caused rows 1-50 to be
# clobbered with stale values and dropped 49 of the 50 generated-image rows).
$rows = @(
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_house_0.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_falls_1.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_boat_2.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_girl_3.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_girl_4.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_house_5.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_house_6.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_boat_7.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_boat_8.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_boat_9.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_girl_10.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_boat_11.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_falls_12.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_house_13.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_girl_14.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_house_15.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_girl_16.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_girl_17.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_18.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_boat_19.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_girl_20.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_house_21.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_falls_22.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_house_23.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_house_24.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_falls_25.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_26.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_house_27.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_boat_28.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_boat_29.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_girl_30.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_girl_31.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_girl_32.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_girl_33.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_house_34.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_house_35.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_secret_cover_boat_36.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_boat_37.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_boat_38.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_girl_39.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_falls_40.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_house.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_house_41.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_42.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_43.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/joly_beans,_holy_greens,_and_billy_jeen_cover_boat_44.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_45.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_boat.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_boat_46.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_falls_47.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_falls.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/woah_oh_ohhhhhh_cover_falls_48.jpg', $TRUE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/cover_girl.jpg', $FALSE),
    @('/home/kdaus/Anti_Forensic_Deep_Learning_Tool/dataGen/gen_data/My_super_duper_secret_cover_girl_49.jpg', $TRUE)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
